# Add a new "OFFSHORE_TRANS" connection row to the Connections sheet,
# right above the existing "NATGAS_TO_CCS" row, and mark every scenario
# column ("Y") for it (matching the pattern used by the BIO_TO_ATM row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Connections")
$ws.Activate()

# Insert a new row above row 15 (currently NATGAS_TO_CCS), shifting
# NATGAS_TO_CCS / BIO_TO_ATM / BIO_TO_CCS down by one row.
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "OFFSHORE_TRANS"
$ws.Range("B15:I15").Value = "Y"

$ws.Range("B15:I15").Select()
